$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows before row 22 first, so every subsequent write lands
# directly on its final address (old rows 22-24 shift down to 24-26,
# leaving a blank 22-23 and blank 27-32 to populate below).
$ws.Rows("22:23").Insert()

# Game Over Screen block (rows 28-32), written before the "Game Win" row so
# new shared-string entries come out in the same order the author typed
# them in.
$ws.Range("D28").Value = "Displays correctly after game loss"
$ws.Range("C28").Value = "Game Over Screen"

$ws.Range("C29").Value = "GO Restart Button Hover"
$ws.Range("C30").Value = "GO Restart Button Press"
$ws.Range("C31").Value = "GO Main Menu Button Hover"
$ws.Range("C32").Value = "GO Main Menu Button Press"

$ws.Range("D29").Value = "Darkens to a slight grey on mouse hover"
$ws.Range("D32").Value = "Darkens to almost black and then Redirects player to main menu"
$ws.Range("D30").Value = "Darkens to almost black and then restarts the snake game from beginning"
$ws.Range("D31").Value = "Darkens to a slight grey on mouse hover"

# Game Win row (row 27)
$ws.Range("C27").Value = "Game Win"
$ws.Range("D27").Value = "Game correctly ends and proceeds to next scene upon reaching 0 pellets"

# Food Pellets / Snake Size rows (rows 22-23)
$ws.Range("C22").Value = "Food Pellets Disappear When Eaten"
$ws.Range("D22").Value = "Food pellets disappear when eaten by snake"

$ws.Range("C23").Value = "Snake Increases in Size"
$ws.Range("D23").Value = "Snake increases in size when pellet eaten"

# Column C needs to widen to fit the new longer text (manually resized, so
# Excel drops the bestFit flag and records a fixed custom width close to
# the old 28.14 + 5 "characters" the new entries need).
$ws.Columns("C").ColumnWidth = 32.28

# Update the worksheet view to match the new scroll / selection position
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("D24").Select()
